# Insert two new rows right after row 174 (pushing existing rows 175+ down by 2)
# and populate them with a new Cilantro price entry pair (Primera / Segunda).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$insertRange = $ws.Range("A175:R176")
$insertRange.Insert()

# New row 175 - "Primera"
$ws.Cells.Item(175, 1).Value = 11
$ws.Cells.Item(175, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(175, 3).Value = "Bíobío"
$ws.Cells.Item(175, 4).Value = 44841
$ws.Cells.Item(175, 5).Value = 8
$ws.Cells.Item(175, 6).Value = 100112040
$ws.Cells.Item(175, 7).Value = "Cilantro"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 200
$ws.Cells.Item(175, 11).Value = 700
$ws.Cells.Item(175, 12).Value = 800
$ws.Cells.Item(175, 13).Value = 750
$ws.Cells.Item(175, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(175, 15).Value = "Región de Ñuble"
$ws.Cells.Item(175, 16).Value = 750
$ws.Cells.Item(175, 17).Value = 1
$ws.Cells.Item(175, 18).Value = "Hortaliza"

# New row 176 - "Segunda"
$ws.Cells.Item(176, 1).Value = 11
$ws.Cells.Item(176, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(176, 3).Value = "Bíobío"
$ws.Cells.Item(176, 4).Value = 44841
$ws.Cells.Item(176, 5).Value = 8
$ws.Cells.Item(176, 6).Value = 100112040
$ws.Cells.Item(176, 7).Value = "Cilantro"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Segunda"
$ws.Cells.Item(176, 10).Value = 100
$ws.Cells.Item(176, 11).Value = 600
$ws.Cells.Item(176, 12).Value = 600
$ws.Cells.Item(176, 13).Value = 600
$ws.Cells.Item(176, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(176, 15).Value = "Región de Ñuble"
$ws.Cells.Item(176, 16).Value = 600
$ws.Cells.Item(176, 17).Value = 1
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D (style s="2").
$ws.Range("D175:D176").NumberFormat = $ws.Range("D177").NumberFormat
